$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B34:D35").Value = $true
